# Extracted products based on commitment level
# Adds a new "Sheet2" with AM/PM product suggestions per commitment level.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create Sheet2 right after Sheet1 -------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Match the outline defaults used on Sheet1 (summary rows/cols above/left)
$ws2.Outline.SummaryRow = 0
$ws2.Outline.SummaryColumn = 0

# --- Column widths -----------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 29.5
$ws2.Columns.Item(2).ColumnWidth = 42.05
$ws2.Columns.Item(3).ColumnWidth = 16.28

# --- Header row ----------------------------------------------------------
$ws2.Range("A1").Value = "Product Suggestion AM"
$ws2.Range("B1").Value = "Product Suggestion PM"
$ws2.Range("C1").Value = "Commitment Level"

# --- Data: fill Commitment Level column first, then AM, then PM ----------
# (mirrors the original authoring order so shared-string ids line up)
$ws2.Range("C2").Value = "Minimal "
$ws2.Range("C3").Value = "Moderate "
$ws2.Range("C4").Value = "Intensive "

$ws2.Range("A2").Value = "•Facewash" + [char]10 + "•Moisturiser" + [char]10 + "•Sunscreen"
$ws2.Range("A3").Value = "•Facewash" + [char]10 + "•Antioxidant serum (like Vitamin C)" + [char]10 + "•Moisturizer" + [char]10 + "•Sunscreen"
$ws2.Range("A4").Value = "•Facewash" + [char]10 + "•Toner (hydrating/exfoliating)" + [char]10 + "•Essence" + [char]10 + "•Antioxidant serum (Vitamin C)" + [char]10 + "•Moisturizer" + [char]10 + "•Sunscreen"

$ws2.Range("B2").Value = "•Facewash" + [char]10 + "•Serum/ Moisturiser" + [char]10
$ws2.Range("B3").Value = "•Facewash" + [char]10 + "•Treatment serum " + [char]10 + "(e.g., niacinamide, hyaluronic acid)" + [char]10 + "•Moisturizer"
$ws2.Range("B4").Value = "•Cleanse" + [char]10 + "•Toner" + [char]10 + "•Active serum (AHA/BHA, retinol, peptides, etc.)" + [char]10 + "•Hydrating serum (hyaluronic acid)" + [char]10 + "•Moisturizer or sleeping mask"

# --- Formatting: reuse Sheet1's existing styles via copy/paste-special ---
# Header row -> bold + wrap (same look as Sheet1's "Skincare Routine
# Commitment Level" / "Preferred Skincare Ingredients/Products" headers)
$ws1.Range("C1").Copy() | Out-Null
$ws2.Range("A1:C1").PasteSpecial(-4122) | Out-Null

# Commitment Level column -> same style as Sheet1's commitment-level cells
$ws1.Range("C2").Copy() | Out-Null
$ws2.Range("C2:C4").PasteSpecial(-4122) | Out-Null

# AM / PM suggestion columns -> wrapped text style like Sheet1's ingredient
# column
$ws1.Range("D2").Copy() | Out-Null
$ws2.Range("A2:B4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- View: freeze header row, make Sheet2 the active tab ------------------
$ws2.Activate()
$ws2.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("B14").Select() | Out-Null
